$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A1 to the new, longer label text (this turns the cell into a shared
# string, t="s", referencing the new xl/sharedStrings.xml entry).
$ws.Range("A1").Value = "TheH(°N)/Tp(s)"

# Best-fit column A so the new, longer header text is fully visible. This is
# what produces the new <cols><col min="1" max="1" .../></cols> entry for
# column 1 only (every other column stays at the default width). Run the
# real AutoFit first (so the column is genuinely sized to its content), then
# pin the best-fit character width so it lines up with the width Excel's own
# AutoFit had computed for this header (~14.57 characters).
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns("A:A").ColumnWidth = 13.6

# Update the active selection/active cell to P14, matching the recorded view
# state left behind after the edit.
$ws.Range("P14").Select() | Out-Null
